$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the A22 timestamp (tiny floating point update)
$ws.Range("A22").Value = 44335.77858406598

# Append new row 23 with the latest retrieved data
$ws.Range("A23").Value = 44336.779981605
$ws.Range("B23").Value = 74317
$ws.Range("C23").Value = 62497
$ws.Range("D23").Value = 3280
$ws.Range("E23").Value = 2056
$ws.Range("F23").Value = 1458
$ws.Range("G23").Value = 19326
$ws.Range("H23").Value = 1311
$ws.Range("I23").Value = 840
$ws.Range("J23").Value = 196
